# Update excel tables
# The "data" worksheet had an extra leading ID column (A: numeric index,
# B: Name, C: Descript). The ID column is no longer needed, so delete the
# entire column A, shifting Name -> A and Descript -> B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Columns.Item(1).Delete()
